# TMTI0056883_VerifyNewJobTypeIsForFVAOnly.xlsx
# "Changes after Recovered data Aug172023"
#
# Adds two more test-data rows to the RecordType and JobType lookup sheets
# (repeating existing CF/FR record types, and a new "CVAS - IP Valuation"
# job type), then leaves the selection/active sheet the way the author
# left it when they saved: RecordType at B5, JobType (the active tab) at C11.

$wb = $excel.ActiveWorkbook

# --- RecordType sheet: duplicate the existing CF / FR rows ---
$wsRecordType = $wb.Worksheets.Item("RecordType")
$wsRecordType.Range("A4").Value = "CF"
$wsRecordType.Range("A5").Value = "FR"
$wsRecordType.Range("B5").Select()

# --- JobType sheet: add the new "CVAS - IP Valuation" job type twice ---
$wsJobType = $wb.Worksheets.Item("JobType")
$wsJobType.Range("A4").Value = "CVAS - IP Valuation"
$wsJobType.Range("A5").Value = "CVAS - IP Valuation"
$wsJobType.Range("C11").Select()
